$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    ,@(2, "Bitcoin", "https://coinranking.com/coin/Qwsogvtv82FCd+bitcoin-btc", "70.843.50", "  +3.24%  ", 0)
    ,@(3, "Ethereum", "https://coinranking.com/coin/razxDUgYGNAdQ+ethereum-eth", "3.573.37", "  +2.48%  ", 0)
    ,@(4, "TetherUSD", "https://coinranking.com/coin/HIVsRcGKkPFtW+tetherusd-usdt", "1.00", "  -0.15%  ", 1)
    ,@(5, "BNB", "https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb", "582.33", "  +2.33%  ", 1)
    ,@(6, "Solana", "https://coinranking.com/coin/zNZHO_Sjf+solana-sol", "186.60", "  +2.38%  ", 1)
    ,@(7, "XRP", "https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp", "0.629", "  +2.56%  ", 1)
    ,@(8, "LidoStakedEther", "https://coinranking.com/coin/VINVMYf0u+lidostakedether-steth", "3.563.73", "  +2.43%  ", 0)
    ,@(9, "USDC", "https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc", "1.00", "  -0.10%  ", 1)
    ,@(10, "Dogecoin", "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge", "0.223", "  +22.40%  ", 1)
    ,@(11, "Cardano", "https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada", "0.652", "  +1.58%  ", 1)
    ,@(12, "Avalanche", "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax", "54.59", "  +1.42%  ", 1)
    ,@(13, "ShibaInu", "https://coinranking.com/coin/xz24e0BjL+shibainu-shib", "0.0000319", "  +6.75%  ", 1)
    ,@(14, "Polkadot", "https://coinranking.com/coin/25W7FG7om+polkadot-dot", "9.47", "  +0.91%  ", 1)
    ,@(15, "WrappedliquidstakedEther2.0", "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth", "4.140.33", "  +2.32%  ", 0)
    ,@(16, "WrappedBTC", "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc", "70.852.56", "  +3.20%  ", 0)
    ,@(17, "Chainlink", "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link", "19.24", "  +0.36%  ", 1)
    ,@(18, "Uniswap", "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni", "12.81", "  +4.95%  ", 1)
    ,@(19, "WrappedEther", "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth", "3.563.14", "  +2.20%  ", 0)
    ,@(20, "BitcoinCash", "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch", "574.44", "  +6.54%  ", 1)
    ,@(21, "TRON", "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx", "0.120", "  +0.63%  ", 1)
    ,@(22, "Polygon", "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic", "1.00", "  -0.48%  ", 1)
    ,@(23, "InternetComputer(DFINITY)", "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp", "17.60", "  -7.68%  ", 1)
    ,@(24, "PancakeSwap", "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake", "4.56", "  +4.33%  ", 1)
    ,@(25, "Toncoin", "https://coinranking.com/coin/67YlI0K1b+toncoin-ton", "4.90", "  -1.60%  ", 1)
    ,@(26, "Litecoin", "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc", "94.03", "  +0.23%  ", 1)
    ,@(27, "RenderToken", "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr", "11.23", "  +4.61%  ", 1)
    ,@(28, "ImmutableX", "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx", "2.95", "  +2.12%  ", 1)
    ,@(29, "Filecoin", "https://coinranking.com/coin/ymQub4fuB+filecoin-fil", "9.22", "  +1.83%  ", 1)
    ,@(30, "EthereumClassic", "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc", "32.42", "  +3.12%  ", 1)
    ,@(31, "NEARProtocol", "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near", "7.20", "  +0.66%  ", 1)
    ,@(32, "Cosmos", "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom", "12.30", "  -1.25%  ", 1)
    ,@(33, "Hedera", "https://coinranking.com/coin/jad286TjB+hedera-hbar", "0.116", "  +2.81%  ", 1)
    ,@(34, "OKB", "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb", "63.11", "  -2.21%  ", 1)
    ,@(35, "Fetch.AI", "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet", "3.41", "  +15.29%  ", 1)
    ,@(36, "dogwifhat", "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif", "3.59", "  +16.82%  ", 1)
    ,@(37, "Bittensor", "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao", "550.04", "  -2.83%  ", 1)
    ,@(38, "TheGraph", "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt", "0.415", "  +5.23%  ", 1)
    ,@(39, "InjectiveProtocol", "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj", "38.12", "  +1.38%  ", 1)
    ,@(40, "PEPE", "https://coinranking.com/coin/03WI8NQPF+pepe-pepe", "0.0₃0804", "  +5.63%  ", 0)
    ,@(41, "Dai", "https://coinranking.com/coin/MoTuySvg7+dai-dai", "0.999", "  -0.01%  ", 1)
    ,@(42, "Maker", "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr", "3.584.52", "  +10.59%  ", 0)
    ,@(43, "Kaspa", "https://coinranking.com/coin/V8GxkwWow+kaspa-kas", "0.138", "  +5.12%  ", 1)
    ,@(44, "Stacks", "https://coinranking.com/coin/mMPrMcB7+stacks-stx", "3.43", "  +3.65%  ", 1)
    ,@(45, "VeChain", "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet", "0.0463", "  +6.46%  ", 1)
    ,@(46, "ApeXProtocol", "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex", "3.49", "  +0.74%  ", 1)
    ,@(47, "ThetaToken", "https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta", "2.93", "  -1.05%  ", 1)
    ,@(48, "THORChain", "https://coinranking.com/coin/ybmU-kKU+thorchain-rune", "9.35", "  +4.62%  ", 1)
    ,@(49, "Stellar", "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm", "0.138", "  +2.64%  ", 1)
    ,@(50, "OceanProtocol", "https://coinranking.com/coin/aAKLSV5-0+oceanprotocol-ocean", "1.50", "  +15.06%  ", 1)
    ,@(51, "FirstDigitalUSD", "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd", "0.998", "  -0.21%  ", 1)
)

foreach ($row in $data) {
    $r = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    if ($row[5] -eq 1) {
        $ws.Cells.Item($r, 4).NumberFormat = "@"
    }
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
}
